$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-29 Saturday", "2024-06-30 Sunday"),
    @("696×2=", "726×5="),
    @("999×8=", "511×5="),
    @("367×5=", "611×4="),
    @("272×3=", "499×9="),
    @("452×6=", "262×2="),
    @("414×2=", "199×4="),
    @("878×6=", "166×2="),
    @("714×2=", "980×6="),
    @("979×6=", "231×6="),
    @("537×7=", "704×4="),
    @("982×3=", "798×8="),
    @("246×8=", "471×9="),
    @("744×9=", "876×2="),
    @("215×2=", "511×6="),
    @("727×9=", "713×6="),
    @("248×7=", "738×9="),
    @("998×7=", "825×2="),
    @("193×7=", "317×4="),
    @("477×5=", "888×8="),
    @("621×8=", "372×7="),
    @("864×3=", "151×6="),
    @("136×3=", "812×3="),
    @("781×7=", "858×6="),
    @("922×4=", "621×6="),
    @("737×7=", "415×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}
